# Generate Report for handback
# Update the "Correspond Handoff Datetime" (column D) and
# "Correspond Handback DateTime" (column G) values for the second data
# row (row 3) on both the "zh-cn" and "de-de" worksheets to reflect the
# newly generated handback report timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-19 06:53:00"
$wsZhCn.Range("G3").Value = "2016-01-19 06:53:44"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-19 06:53:10"
$wsDeDe.Range("G3").Value = "2016-01-19 06:54:01"
